# Week 13 logging update
# Appends this week's per-play/drive tracking numbers to the running logs on
# the YDS and ST sheets, and updates the corresponding weekly summary totals
# on the OFF, DEF, ST, TURNS and PEN sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet - running space separated logs of play-by-play yardage
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value() + " 3 7 1 3 0 17 7 3 15 -1 6 2 1 1 10 1"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value() + " 22 1 3 1 16 2 6 2 2 4 1 9 5 -4 -3 13 4 1 0 -6 7 2 3 8 1 6 0 4 3"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value() + " 1 8 9 3 10 5 3 9 7 5 0 18 10 -5 5 2 15 9 34 7 3 28 14 13 16 2 11 7"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value() + " 23 -1 13 6 7 4 7 6 5 6 5 6 0 10 35 10 3 4 10 9 10 12 6"

# ---------------------------------------------------------------------------
# OFF sheet - weekly offensive totals
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 122
$wsOFF.Range("D2").Value = 14
$wsOFF.Range("E2").Value = 5
$wsOFF.Range("F2").Value = 34
$wsOFF.Range("G2").Value = 46
$wsOFF.Range("L2").Value = 266
$wsOFF.Range("M2").Value = 185
$wsOFF.Range("Q2").Value = 436

$wsOFF.Range("C3").Value = 162
$wsOFF.Range("E3").Value = 26
$wsOFF.Range("F3").Value = 84
$wsOFF.Range("G3").Value = 35
$wsOFF.Range("H3").Value = 19
$wsOFF.Range("I3").Value = 40
$wsOFF.Range("J3").Value = 39
$wsOFF.Range("N3").Value = 15

# ---------------------------------------------------------------------------
# DEF sheet - weekly defensive totals
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value = 155
$wsDEF.Range("F2").Value = 42
$wsDEF.Range("G2").Value = 60
$wsDEF.Range("H2").Value = 2
$wsDEF.Range("I2").Value = 6
$wsDEF.Range("J2").Value = 27
$wsDEF.Range("L2").Value = 241
$wsDEF.Range("M2").Value = 160
$wsDEF.Range("Q2").Value = 478

$wsDEF.Range("B3").Value = 7
$wsDEF.Range("C3").Value = 139
$wsDEF.Range("E3").Value = 21
$wsDEF.Range("F3").Value = 87
$wsDEF.Range("H3").Value = 18
$wsDEF.Range("I3").Value = 51
$wsDEF.Range("J3").Value = 40
$wsDEF.Range("N3").Value = 13

# ---------------------------------------------------------------------------
# ST sheet - weekly special teams totals + running logs
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 62
$wsST.Range("D2").Value = 48
$wsST.Range("H2").Value = 2
$wsST.Range("J2").Value = 64
$wsST.Range("K2").Value = 61
$wsST.Range("N2").Value = 15
$wsST.Range("O2").Value = 12

$wsST.Range("B3").Value = 37
$wsST.Range("D3").Value = $wsST.Range("D3").Value() + " 42 67 38 55"

$wsST.Range("D4").Value = $wsST.Range("D4").Value() + " 0 16 0 13"
$wsST.Range("D5").Value = $wsST.Range("D5").Value() + " 6 0 0 0"
$wsST.Range("B6").Value = $wsST.Range("B6").Value() + " 19 24"

# ---------------------------------------------------------------------------
# TURNS sheet - weekly turnover totals
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("C2").Value = 1
$wsTURNS.Range("D3").Value = 4

# ---------------------------------------------------------------------------
# PEN sheet - weekly penalty totals
# ---------------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("D2").Value = 10
$wsPEN.Range("D3").Value = 5
